$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Query($cellRef) {
    $val = $ws.Range($cellRef).Value()
    $val = $val.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $val = $val.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $val = $val.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $val = $val.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $val = $val.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $val = $val.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    $ws.Range($cellRef).Value = $val
}

# StudiesTab: StatQuery (C2) and TabQuery (B2)
Update-Query "C2"
Update-Query "B2"

# ParticipantsTab TabQuery
Update-Query "B3"

# DiagnosisTab TabQuery
Update-Query "B4"

# TreatmentTab TabQuery
Update-Query "B5"

# TreatmentRespTab TabQuery
Update-Query "B6"

# SurvivalTab TabQuery
Update-Query "B7"

# Column C width change (no longer best-fit, explicit width 67.5 in saved XML)
$ws.Columns.Item(3).ColumnWidth = 66.66666666666667
